# edit.ps1 -- applies the "Added deserialization, read, dual-mode filtering" change
# to the TestModels sheet: a new data row, refreshed date-time formatting/fonts,
# a frozen header pane, an updated selection, and a portrait page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestModels")
$ws.Activate()

# --- 1. Header row (row 1): italic font, keep B/C as the custom date-time format ---
$ws.Range("A1:E1").Font.Italic = $true
$ws.Range("B1:C1").NumberFormat = "m/d/yy h:mm;@"

# --- 2. Existing data row 2: refresh the timestamp number format + precise values ---
$ws.Range("B2:C2").NumberFormat = "m/d/yy h:mm;@"
$ws.Cells.Item(2, 2).Value = 43047.19836610485
$ws.Cells.Item(2, 3).Value = 43047.19836805556

# --- 3. New row 3 -- second TestModels record ("Rob") ---
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 43058.05815972222
$ws.Cells.Item(3, 3).Value = 43058.2549537037
$ws.Cells.Item(3, 4).Value = "Rob"
$ws.Cells.Item(3, 5).Value = 34
$ws.Range("B3:C3").NumberFormat = "m/d/yy h:mm;@"

# --- 4. Column widths (best effort -- COM only exposes whole-pixel granularity) ---
$ws.Columns.Item(1).ColumnWidth = 3.28
$ws.Columns.Item(2).ColumnWidth = 11.72
$ws.Columns.Item(3).ColumnWidth = 11.72
$ws.Columns.Item(4).ColumnWidth = 5.39
$ws.Columns.Item(5).ColumnWidth = 3.83

# --- 5. Freeze the header row, and leave the final selection on C5 ---
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$ws.Range("C5").Select()

# --- 6. Portrait page setup ---
$ws.PageSetup.Orientation = 1

# --- 7. Outline defaults (summary rows below / summary columns to the right) on every sheet ---
foreach ($sheet in $wb.Worksheets) {
    $sheet.Outline.SummaryRow = 1
    $sheet.Outline.SummaryColumn = 1
}

Write-Host "TestModels updated: new row for Rob, refreshed formatting, frozen header pane."
